# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values that look like plain numbers
# are stored as text in this sheet (t="inlineStr"), matching the site's
# 'thousands-dot' price formatting (e.g. '52.306.61'), so we force them to
# stay text with a leading apostrophe and then restore the cell's original
# style afterwards so no stray number formatting / quote-prefix is left on
# the cell.
$updates = [ordered]@{
    "D2" = "52.306.61"
    "E2" = "  +1.51%  "
    "D3" = "2.791.08"
    "E3" = "  +1.45%  "
    "E4" = "  +0.03%  "
    "D5" = "'346.41"
    "E5" = "  +4.00%  "
    "D6" = "'116.02"
    "E6" = "  +0.13%  "
    "E7" = "  +3.14%  "
    "E8" = "  -0.06%  "
    "D9" = "'0.592"
    "E9" = "  +2.87%  "
    "D10" = "'42.62"
    "E10" = "  +2.48%  "
    "D11" = "'0.0857"
    "E11" = "  +3.45%  "
    "D12" = "'20.02"
    "E12" = "  -1.27%  "
    "E13" = "  +1.75%  "
    "E14" = "  +2.83%  "
    "D15" = "3.240.64"
    "E15" = "  +1.66%  "
    "D16" = "2.803.21"
    "E16" = "  +2.15%  "
    "D17" = "'0.891"
    "E17" = "  +0.25%  "
    "D18" = "52.177.22"
    "E18" = "  +1.27%  "
    "D19" = "'3.17"
    "E19" = "  +5.99%  "
    "D20" = "'7.23"
    "E20" = "  +5.04%  "
    "D21" = "'13.44"
    "E21" = "  -2.96%  "
    "D22" = "0.0₃0979"
    "E22" = "  +1.66%  "
    "D23" = "'269.57"
    "E23" = "  -3.62%  "
    "D24" = "'69.96"
    "E24" = "  -0.37%  "
    "D25" = "'2.75"
    "E25" = "  +4.05%  "
    "D26" = "'26.80"
    "E26" = "  -0.81%  "
    "D27" = "'0.999"
    "E27" = "  -0.12%  "
    "D28" = "'10.20"
    "E28" = "  -1.88%  "
    "E29" = "  +0.68%  "
    "D30" = "'0.141"
    "E30" = "  -0.06%  "
    "B31" = "VeChain"
    "C31" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D31" = "'0.0478"
    "E31" = "  +36.64%  "
    "B32" = "InjectiveProtocol"
    "C32" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D32" = "'34.82"
    "E32" = "  -3.01%  "
    "B33" = "OKB"
    "C33" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D33" = "'49.97"
    "E33" = "  -1.16%  "
    "D34" = "'5.72"
    "E34" = "  +1.36%  "
    "D35" = "'0.0827"
    "E35" = "  -0.26%  "
    "E36" = "  +0.05%  "
    "E37" = "  -0.18%  "
    "D38" = "'4.97"
    "E38" = "  -0.98%  "
    "D39" = "'18.60"
    "E39" = "  -4.66%  "
    "D40" = "'3.22"
    "E40" = "  -0.51%  "
    "D41" = "'2.61"
    "E41" = "  +10.25%  "
    "D42" = "'127.38"
    "E42" = "  -2.12%  "
    "B43" = "Stellar"
    "C43" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D43" = "'0.115"
    "E43" = "  +1.61%  "
    "B44" = "EnergySwap"
    "C44" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D44" = "'23.14"
    "E44" = "  -1.53%  "
    "D45" = "'2.30"
    "E45" = "  -0.08%  "
    "D46" = "'3.32"
    "E46" = "  -2.38%  "
    "D47" = "2.061.95"
    "E47" = "  -2.55%  "
    "E48" = "  +2.79%  "
    "D49" = "'0.947"
    "E49" = "  +11.01%  "
    "D50" = "'5.55"
    "E50" = "  -0.81%  "
    "D51" = "'8.94"
    "E51" = "  -0.98%  "
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $range = $ws.Range($addr)
    $originalStyle = $range.Style
    $range.Value = $value
    $range.Style = $originalStyle
}
